# Maintenance Plan.docx edit
#
# The red placeholder block ("Maintenance plan (500-700 words - 5%)" plus the
# six bullet/question paragraphs that followed it, all styled in red) is
# replaced by the team's actual maintenance-plan write-up. The write-up lives
# in a single paragraph; the paragraph mark keeps the original red color
# (inherited from the placeholder paragraph) but the visible text runs are
# left uncolored, and a couple of proofing marks (gramStart/gramEnd around
# "free") are preserved from the authored content.

$d = $word.ActiveDocument

# Find the start of the placeholder block (the first red paragraph).
$startRng = $d.Content
$foundStart = $startRng.Find.Execute("Maintenance plan (500-700 words - 5%)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundStart) {
    throw "Could not find the start of the maintenance-plan placeholder block"
}
$startPara = $startRng.Paragraphs(1)

# Find the end of the placeholder block (the last red paragraph, "...XBox live)").
$endRng = $d.Content
$foundEnd = $endRng.Find.Execute("live)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundEnd) {
    throw "Could not find the end of the maintenance-plan placeholder block"
}
$endPara = $endRng.Paragraphs(1)

# Whole block spanning every red placeholder paragraph (start of the first
# through the paragraph mark of the last).
$block = $d.Range($startPara.Range.Start, $endPara.Range.End)

# Replace the whole block in one shot with a single new paragraph, expressed
# as a WordprocessingML package fragment so formatting/run-splits/proofErr
# marks come out exactly as authored.
$newXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">Our idea for this project is to have an app that helps students in need of resources for their mental health. This app has the essential software and styling to help with a person’s need. There is a survey with questions that are asked and once answered the app gathers information and collects a history on the answers. This way allowing the user to reflect on their feelings. The app also gives advice and motivational quotes. For this project we did not buy any program to run our app because it is just a project but if we were to make useable we would have to consider how much it would cost to have it up and running. We would have to hire someone to maintain our app and the app placement itself. That way our app can be used by any college student that would like to access it. Our app is planned to be </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>free</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> but we would still have to pay the app store around $99 a year. The software for this code is not heavy and should be of light maint</w:t></w:r><w:r><w:t>e</w:t></w:r><w:r><w:t>n</w:t></w:r><w:r><w:t>an</w:t></w:r><w:r><w:t xml:space="preserve">ce. </w:t></w:r><w:r><w:t>Through time if we would like to scale the app for more colleges and mass use. We would have to make improvements, and more general utilization features. We would also have to add additional useful elements to elevate the app. To do so we would have to hire software developers and graphic designers. I think it would be enough to do $1,500 a month for six months to 3 or 5 developers and graphic designers. This way a couple can work on the front end of the app, one on the back, and a designer. This way the app will keep its proficiency. There are other factors that may also be included in the price of the app. Such as a marketing team to expand ideas for new features or improvements on the app that should be made. Right now we are running the app in our browser and if we were to make a website domain for people to be able to reach their accounts on their desktop, laptop, or browser it could cost depending on the place we decide around $40 a year. For right now it seems that running our app for a year does not run us out of too much money but that all depends on</w:t></w:r><w:r><w:t xml:space="preserve"> how quickly or when we would like to expand our product. It could also be less money if we decide to do our own work and take care of the app ourselves. This way all we would need to pay is the annual fees for the domain and servers, the annual fees for distribution platform and we would save ourselves the money for hiring developers and or graphic designers. However, depending on how everything goes it may be more beneficial to hire help. How much we would like to pay the developers can also be less or more depending on who we hire, their rates, and how much work we would like done. Also having someone that can take care of the app would be the most expensive because we would have to pay them for the year to ensure that the app is running smoothly and everything is taken care off, and the users are having an enjoyable experience. This could cost around $25,00+. Again, depending on who we hire, experience level, and workload. All of the estimates could be less or more depending on the quality we desire for the app.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

[void]$block.InsertXML($newXml)
